$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain parseable number must be forced to Text
# format first (matching the source data, which stores every Price/Volume
# entry as text) so Excel does not silently convert them to numeric cells.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated values
$ws.Range("D2").Value = "27.928.34"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.629.92"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "211.87"
$ws.Range("E5").Value = "  -0.77%  "
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -1.01%  "
$ws.Range("D9").Value = "0.257"
$ws.Range("E9").Value = "  -1.71%  "
$ws.Range("E10").Value = "  -0.26%  "
$ws.Range("D11").Value = "0.0880"
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("D12").Value = "1.860.51"
$ws.Range("D13").Value = "1.626.71"
$ws.Range("E13").Value = "  -0.76%  "
$ws.Range("D14").Value = "4.04"
$ws.Range("E14").Value = "  -1.48%  "
$ws.Range("E15").Value = "  -2.36%  "
$ws.Range("D16").Value = "65.60"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").Value = "27.917.24"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "230.45"
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D22").Value = "4.35"
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("E23").Value = "  -4.61%  "
$ws.Range("D25").Value = "154.78"
$ws.Range("E25").Value = "  +2.04%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").Value = "15.55"
$ws.Range("E28").Value = "  -1.20%  "
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("E30").Value = "  -0.74%  "
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("E32").Value = "  +2.36%  "
$ws.Range("E33").Value = "  -0.67%  "
$ws.Range("D34").Value = "1.402.33"
$ws.Range("E34").Value = "  -0.47%  "
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("E36").Value = "  +11.04%  "
$ws.Range("D37").Value = "2.35"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("E38").Value = "  +1.98%  "
$ws.Range("D39").Value = "0.556"
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("E40").Value = "  -2.99%  "
$ws.Range("E41").Value = "  -0.34%  "
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("D47").Value = "1.770.58"
$ws.Range("D48").Value = "88.10"
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0103"
$ws.Range("E49").Value = "  -2.56%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.101"
$ws.Range("E50").Value = "  +0.46%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.0504"
$ws.Range("E51").Value = "  -0.23%  "
